# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# (and the overview handback timestamp) cells to reflect the latest run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (row 2 corresponds to the
# 1405fd1f-... record, shared with the de-de sheet's handoff datetime).
$wsOverview.Range("G2").Value = "2016-08-26 13:06:43"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-26 13:06:38"
$wsZhCn.Range("K2").Value = "2016-08-26 13:06:56"

# de-de sheet: Correspond Handoff Datetime stays tied to the Overview value,
# only the Correspond Handback DateTime changes.
$wsDeDe.Range("H2").Value = "2016-08-26 13:06:43"
$wsDeDe.Range("K2").Value = "2016-08-26 13:07:10"
